$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
# L20: PIEDRA SINTERIZADA sale for this client, 0 -> 582.53
$ws1.Range("L20").Value = 582.53
# L30: count text "1 de 28" -> "2 de 28"
$ws1.Range("L30").Value = "2 de 28"

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
# F20: junio sale for this client, 0 -> 582.53
$ws2.Range("F20").Value = 582.53
# F30: junio total, 1687.38 -> 2269.91
$ws2.Range("F30").Value = 2269.91

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Row 15: PIEDRA SINTERIZADA group compliance
$ws3.Range("D15").Value = 867.65
$ws3.Range("E15").Value = -340.62
$ws3.Range("F15").Value = 1.646300969584274

# Row 19: TOTAL row
$ws3.Range("D19").Value = 2264.15
$ws3.Range("E19").Value = 27273.64107555787
$ws3.Range("F19").Value = 0.07665265131736791
